$wb = $excel.ActiveWorkbook

# ALC!row54
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# ALC!row88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 920.1429000000001
$ws.Range("I88").Value = 1722
$ws.Range("J88").Value = 599.4
$ws.Range("K88").Value = 1722
$ws.Range("L88").Value = 599.4
$ws.Range("M88").Value = -1316
$ws.Range("N88").Value = -1411.4

# ALC!row91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 920.1429000000001
$ws.Range("I91").Value = 1722
$ws.Range("J91").Value = 599.4
$ws.Range("K91").Value = 1722
$ws.Range("L91").Value = 599.4
$ws.Range("M91").Value = -318
$ws.Range("N91").Value = -3407.4

# ALC!row99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 7628.5713
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 7628.5713
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 22885.7139
$ws.Range("N99").Value = -25881.7139
$ws.Range("M99").ClearContents()

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1543.4286
$ws.Range("I100").Value = 1550.6666
$ws.Range("K100").Value = 1550.6666
$ws.Range("M100").Value = -1009.6666

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1743.238
$ws.Range("I137").Value = 844
$ws.Range("K137").Value = 2532
$ws.Range("M137").Value = 18

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2111.5
$ws.Range("I61").Value = 1596.5454
$ws.Range("K61").Value = 1596.5454
$ws.Range("M61").Value = -1384.5454

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4007.7778
$ws.Range("I74").Value = 3206.5715
$ws.Range("K74").Value = 3206.5715
$ws.Range("M74").Value = -2332.5715

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4007.7778
$ws.Range("I77").Value = 3206.5715
$ws.Range("K77").Value = 16032.8575
$ws.Range("M77").Value = -11664.8575

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1420.75
$ws.Range("I132").Value = 1420.75
$ws.Range("K132").Value = 4262.25
$ws.Range("M132").Value = -1732.25

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2111.5
$ws.Range("I136").Value = 1596.5454
$ws.Range("K136").Value = 4789.6362
$ws.Range("M136").Value = -2239.6362

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1919.6
$ws.Range("I20").Value = 2800
$ws.Range("J20").Value = 1332.6666
$ws.Range("K20").Value = 2800
$ws.Range("L20").Value = 1332.6666
$ws.Range("M20").Value = -2553
$ws.Range("N20").Value = -1826.6666

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4448.9
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 5498.4287
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 5498.4287
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -7744.4287

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4448.9
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 5498.4287
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 27492.1435
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -38724.14350000001

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1282.1666
$ws.Range("I94").Value = 1282.1666
$ws.Range("K94").Value = 1282.1666
$ws.Range("M94").Value = -831.1666

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2177.2727
$ws.Range("I134").Value = 1443.3334
$ws.Range("K134").Value = 4330.0002
$ws.Range("M134").Value = -1795.0002

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2291
$ws.Range("I16").Value = 2291
$ws.Range("K16").Value = 2291
$ws.Range("M16").Value = -2004

# CRP!row22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1163.2632
$ws.Range("J22").Value = 1247
$ws.Range("L22").Value = 1247
$ws.Range("N22").Value = -1947

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2890.5334
$ws.Range("I58").Value = 1452.4
$ws.Range("J58").Value = 5766.8
$ws.Range("K58").Value = 1452.4
$ws.Range("L58").Value = 5766.8
$ws.Range("M58").Value = -1249.4
$ws.Range("N58").Value = -6172.8

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3016
$ws.Range("I62").Value = 2972.8
$ws.Range("K62").Value = 2972.8
$ws.Range("M62").Value = -2348.8

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3016
$ws.Range("I65").Value = 2972.8
$ws.Range("K65").Value = 14864
$ws.Range("M65").Value = -11744

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2291
$ws.Range("I113").Value = 2291
$ws.Range("K113").Value = 2291
$ws.Range("M113").Value = -121

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3744.0667
$ws.Range("I134").Value = 2421.182
$ws.Range("J134").Value = 7382
$ws.Range("K134").Value = 7263.545999999999
$ws.Range("L134").Value = 22146
$ws.Range("M134").Value = -4728.545999999999
$ws.Range("N134").Value = -27216

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2890.5334
$ws.Range("I136").Value = 1452.4
$ws.Range("J136").Value = 5766.8
$ws.Range("K136").Value = 4357.200000000001
$ws.Range("L136").Value = 17300.4
$ws.Range("M136").Value = -1807.200000000001
$ws.Range("N136").Value = -22400.4

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5332.3335
$ws.Range("I70").Value = 5332.3335
$ws.Range("K70").Value = 5332.3335
$ws.Range("M70").Value = -5062.3335

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5332.3335
$ws.Range("I73").Value = 5332.3335
$ws.Range("K73").Value = 5332.3335
$ws.Range("M73").Value = -4396.3335

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 954.6667
$ws.Range("I80").Value = 810.5
$ws.Range("J80").Value = 1098.8334
$ws.Range("K80").Value = 810.5
$ws.Range("L80").Value = 1098.8334
$ws.Range("M80").Value = 187.5
$ws.Range("N80").Value = -3094.8334

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 954.6667
$ws.Range("I83").Value = 810.5
$ws.Range("J83").Value = 1098.8334
$ws.Range("K83").Value = 4052.5
$ws.Range("L83").Value = 5494.166999999999
$ws.Range("M83").Value = 939.5
$ws.Range("N83").Value = -15478.167

# GSM!row92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9599.200000000001
$ws.Range("J92").Value = 9599.200000000001
$ws.Range("L92").Value = 9599.200000000001
$ws.Range("N92").Value = -13343.2

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 202284.8
$ws.Range("I132").Value = 202284.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 606854.3999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -604324.3999999999
$ws.Range("N132").ClearContents()

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2294.4443
$ws.Range("I68").Value = 1235.8572
$ws.Range("K68").Value = 1235.8572
$ws.Range("M68").Value = -486.8571999999999

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2294.4443
$ws.Range("I71").Value = 1235.8572
$ws.Range("K71").Value = 6179.286
$ws.Range("M71").Value = -2435.286

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3101.4285
$ws.Range("J82").Value = 4158.7
$ws.Range("L82").Value = 4158.7
$ws.Range("N82").Value = -4880.7

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3101.4285
$ws.Range("J85").Value = 4158.7
$ws.Range("L85").Value = 4158.7
$ws.Range("N85").Value = -6654.7

# LTW!row93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1133.9048
$ws.Range("I93").Value = 1158
$ws.Range("K93").Value = 1158
$ws.Range("M93").Value = 90

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1800.25
$ws.Range("I122").Value = 1800.25
$ws.Range("K122").Value = 5400.75
$ws.Range("M122").Value = -2950.75

# WVR!row96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1149.25
$ws.Range("J96").Value = 1200
$ws.Range("L96").Value = 1200
$ws.Range("N96").Value = -3946

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1807.8
$ws.Range("I122").Value = 1645.3334
$ws.Range("J122").Value = 1877.4286
$ws.Range("K122").Value = 4936.0002
$ws.Range("L122").Value = 5632.2858
$ws.Range("M122").Value = -2486.0002
$ws.Range("N122").Value = -10532.2858

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2496.2
$ws.Range("I132").Value = 2183
$ws.Range("K132").Value = 6549
$ws.Range("M132").Value = -4019
